# Insert the new 'Knärot - ekologi samt krav på livsmiljön' section (a heading,
# five body paragraphs, a 'Referenser - knärot' sub-heading and six reference
# entries) right before the document's closing section break, and bump the
# dated cover line in the first-page header from 2023-09-13 to 2023-09-15.
#
# New paragraphs are created in two passes: pass 1 inserts every new (still
# empty) paragraph mark and stamps its style; pass 2 fills each one in with its
# runs (text + italics). Doing both in a single pass would let a trailing italic
# run of a just-filled paragraph bleed its character formatting into the *next*
# paragraph mark created by InsertParagraphAfter.

$d = $word.ActiveDocument

# Style (Word style NameLocal) to apply to each of the new paragraphs, in order.
$paraStyles = @(
    "Heading 1",
    "Normal",
    "Normal",
    "Normal",
    "Normal",
    "Normal",
    "Heading 2",
    "Normal",
    "Normal",
    "Normal",
    "Normal",
    "Normal",
    "Normal",
)

# --- Pass 1: create the new (empty) paragraphs with their target styles ---
$newParas = New-Object System.Collections.ArrayList
$lastPara = $d.Paragraphs.Last
foreach ($styleName in $paraStyles) {
    $lastPara.Range.InsertParagraphAfter()
    $lastPara = $d.Paragraphs.Last
    $lastPara.Style = $styleName
    [void]$newParas.Add($lastPara)
}

# --- Pass 2: fill in each paragraph's runs (text + italics) ---

function Add-Run {
    param($doc, $pos, $text, $italic)
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    $newPos = $pos + $text.Length
    if ($italic) {
        $fr = $doc.Range($pos, $newPos)
        $fr.Font.Italic = $true
    }
    return $newPos
}

# Paragraph 1/13 (style=Heading1)
$p = $newParas[0]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Knärot – ekologi samt krav på livsmiljön' $false

# Paragraph 2/13 (style=Normal)
$p = $newParas[1]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).' $false

# Paragraph 3/13 (style=Normal)
$p = $newParas[2]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Samuel Johnsons doktorsavhandling ' $false
$pos = Add-Run $d $pos '“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“' $true
$pos = Add-Run $d $pos ' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ' $false
$pos = Add-Run $d $pos '“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ' $true
$pos = Add-Run $d $pos 'Vidare ' $false
$pos = Add-Run $d $pos '“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”' $true

# Paragraph 4/13 (style=Normal)
$p = $newParas[3]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ' $false
$pos = Add-Run $d $pos '“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”' $true

# Paragraph 5/13 (style=Normal)
$p = $newParas[4]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).' $false

# Paragraph 6/13 (style=Normal)
$p = $newParas[5]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).' $false

# Paragraph 7/13 (style=Heading2)
$p = $newParas[6]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Referenser - knärot' $false

# Paragraph 8/13 (style=Normal)
$p = $newParas[7]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'de Graaf M & Roberts M.R., 2009. ' $false
$pos = Add-Run $d $pos 'Short-term response of the herbaceous layer within leave patches after harvest. ' $true
$pos = Add-Run $d $pos 'Forest Ecology and Management 257, 1014-1025' $false

# Paragraph 9/13 (style=Normal)
$p = $newParas[8]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ' $false
$pos = Add-Run $d $pos 'Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ' $true
$pos = Add-Run $d $pos 'Ecological Applications, 22, 2049-2064 ' $false

# Paragraph 10/13 (style=Normal)
$p = $newParas[9]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ' $false
$pos = Add-Run $d $pos 'Interactive effects of drought and edge exposure on old-growth forest understory species. ' $true
$pos = Add-Run $d $pos 'Landscape Ecology, 37, sid 1839-1853' $false

# Paragraph 11/13 (style=Normal)
$p = $newParas[10]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ' $false
$pos = Add-Run $d $pos 'Biological legacies buffer local species extinction after logging. ' $true
$pos = Add-Run $d $pos 'Journal of Applied Ecology. 51, 53-62.' $false

# Paragraph 12/13 (style=Normal)
$p = $newParas[11]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'Skogsstyrelsen, 2022. ' $false
$pos = Add-Run $d $pos 'Vägledning för hänsyn till knärot. ' $true
$pos = Add-Run $d $pos 'https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/' $false

# Paragraph 13/13 (style=Normal)
$p = $newParas[12]
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'SLU Artdatabanken, 2021. ' $false
$pos = Add-Run $d $pos 'Artfaktablad. Naturvård – artfakta. ' $true
$pos = Add-Run $d $pos 'SLU Artdatabanken, Uppsala ' $false

# Update the dated cover line in the first-page header (Find/Execute walks
# every section's headers so it is robust even if the section count changes).
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
